# Insert a new weekly record at row 90 of "Sheet1", pushing the existing
# rows 90-138 down to 91-139 (dimension grows from A1:T138 to A1:T139).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(90).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(90, 1).Value  = 8
$ws.Cells.Item(90, 2).Value  = 'Terminal La Palmera de La Serena'
$ws.Cells.Item(90, 3).Value  = 'Coquimbo'
$ws.Cells.Item(90, 4).Value  = 44960
$ws.Cells.Item(90, 5).Value  = 4
$ws.Cells.Item(90, 6).Value  = 'Fruta'
$ws.Cells.Item(90, 7).Value  = 100109
$ws.Cells.Item(90, 8).Value  = 'Uva'
$ws.Cells.Item(90, 9).Value  = 100109001
$ws.Cells.Item(90, 10).Value = 'Uva'
$ws.Cells.Item(90, 11).Value = 'Flame Seedless'
$ws.Cells.Item(90, 12).Value = 'Primera'
$ws.Cells.Item(90, 13).Value = 600
$ws.Cells.Item(90, 14).Value = 6500
$ws.Cells.Item(90, 15).Value = 7000
$ws.Cells.Item(90, 16).Value = 6750
$ws.Cells.Item(90, 17).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(90, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(90, 19).Value = 375
$ws.Cells.Item(90, 20).Value = 18
